# Scheduled-runner update: refresh cached Universalis market-price snapshots
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ) and the
# derived Leve price/profit columns (LevePriceNQ, LevePriceHQ, LeveProfitNQ,
# LeveProfitHQ) on each job sheet. Values are static snapshots (no formulas in
# this workbook), so the runner overwrites the affected cells directly.

$wb = $excel.ActiveWorkbook

# A sentinel distinguishing "clear this cell" from "leave this cell alone"
# ($null below means "column untouched by this update").
$CLEAR = "##CLEAR##"

# Each entry: SheetName, Row, H, I, J, K, L, M, N
# ($null => cell left as-is; $CLEAR => cell is emptied/removed)
$edits = @(
    ,@("ALC", 11, 38.875, 38.875, 0, 38.875, 0, 101.125, $null)
    ,@("ALC", 32, 3646.5334, 1564.3334, 4167.0835, 1564.3334, 4167.0835, -1238.3334, -4819.0835)
    ,@("ALC", 62, 4174.8, 4853.25, 3399.4285, 4853.25, 3399.4285, -4229.25, -4647.4285)
    ,@("ALC", 65, 4174.8, 4853.25, 3399.4285, 24266.25, 16997.1425, -21146.25, -23237.1425)
    ,@("ALC", 111, 3974, 1965.3334, 10000, 5896.0002, 30000, -2829.0002, -36134)
    ,@("ALC", 129, 1649, 1594.25, 1722, 4782.75, 5166, 217.25, -15166)
    ,@("ALC", 137, 2121.111, 1431.6666, 3500, 4294.9998, 10500, -1744.9998, -15600)
    ,@("ALC", 138, 3167.8462, 2074.1765, 4012.9546, 6222.529500000001, 12038.8638, -1082.529500000001, -22318.8638)
    ,@("ALC", 141, 5566.472, 1883.7097, 28399.6, 5651.1291, 85198.79999999999, -471.1291000000001, -95558.79999999999)
    ,@("ARM", 2, 2400.4285, 1702.7858, 3795.7144, 1702.7858, 3795.7144, -1589.7858, -4021.7144)
    ,@("ARM", 32, 2858.5283, 1640.6522, 10861.714, 1640.6522, 10861.714, -1353.6522, -11435.714)
    ,@("ARM", 61, 4419.778, 2131.1667, 8997, 2131.1667, 8997, -1919.1667, -9421)
    ,@("ARM", 74, 1712.2373, 1677.4043, 1848.6666, 1677.4043, 1848.6666, -803.4042999999999, -3596.6666)
    ,@("ARM", 77, 1712.2373, 1677.4043, 1848.6666, 8387.021499999999, 9243.333000000001, -4019.021499999999, -17979.333)
    ,@("ARM", 110, 1442.6522, 1426.619, 1611, 1426.619, 1611, 618.3810000000001, -5701)
    ,@("ARM", 116, 2400.4285, 1702.7858, 3795.7144, 1702.7858, 3795.7144, 591.2141999999999, -8383.714400000001)
    ,@("ARM", 133, 0, 0, 0, 0, $CLEAR, $null, 0)
    ,@("ARM", 136, 4419.778, 2131.1667, 8997, 6393.500100000001, 26991, -3843.500100000001, -32091)
    ,@("BSM", 3, 2400.4285, 1702.7858, 3795.7144, 1702.7858, 3795.7144, -1588.7858, -4023.7144)
    ,@("BSM", 94, 2257.6858, 892.53845, 6201.4443, 892.53845, 6201.4443, -441.53845, -7103.4443)
    ,@("CRP", 6, 7504960, 7504960, 0, 7504960, 0, -7504847, $null)
    ,@("CRP", 31, 4591.1113, 2572.9285, 6764.5386, 2572.9285, 6764.5386, -2277.9285, -7354.5386)
    ,@("CRP", 34, 4591.1113, 2572.9285, 6764.5386, 2572.9285, 6764.5386, -2370.9285, -7168.5386)
    ,@("CRP", 58, 2525.5, 2441.1, 2666.1667, 2441.1, 2666.1667, -2238.1, -3072.1667)
    ,@("CRP", 134, 2747.9565, 2577.7222, 3360.8, 7733.1666, 10082.4, -5198.1666, -15152.4)
    ,@("CRP", 136, 2525.5, 2441.1, 2666.1667, 7323.299999999999, 7998.500100000001, -4773.299999999999, -13098.5001)
    ,@("CRP", 138, 0, 0, 0, 0, $CLEAR, $null, 0)
    ,@("CUL", 17, 2730.6155, 2750, 2498, 8250, 7494, -8081, -7832)
    ,@("CUL", 70, 4507.3335, 4756, 4010, 14268, 12030, -13953, -12660)
    ,@("CUL", 73, 4507.3335, 4756, 4010, 14268, 12030, -13176, -14214)
    ,@("CUL", 75, 799, 797.6667, 799.26666, 2393.0001, 2397.79998, -1395.0001, -4393.79998)
    ,@("CUL", 78, 799, 797.6667, 799.26666, 7179.0003, 7193.39994, -2187.0003, -17177.39994)
    ,@("CUL", 126, 15000, 0, 15000, 0, 45000, $null, -54880)
    ,@("CUL", 140, 1518.5834, 1518.5834, 0, 4555.7502, 0, $CLEAR, 624.2497999999996)
    ,@("GSM", 97, 562.0714, 548.2222, 587, 548.2222, 587, -52.22220000000004, -1579)
    ,@("GSM", 113, 3377.8333, 3073.8, 4898, 3073.8, 4898, -903.8000000000002, -9238)
    ,@("GSM", 132, 5686.5, 5789.2, 5515.3335, 17367.6, 16546.0005, -14837.6, -21606.0005)
    ,@("LTW", 16, 2733.2778, 2220.7, 3374, 2220.7, 3374, -2050.7, -3714)
    ,@("LTW", 45, 20000, 20000, 0, 20000, 0, -19593, $null)
    ,@("LTW", 46, 951.375, 901, 958.5714, 901, 958.5714, -713, -1334.5714)
    ,@("LTW", 50, 31075.6, 40000, 30084, 40000, 30084, -39363, -31358)
    ,@("LTW", 122, 5430.8, 6898.6, 3963, 20695.8, 11889, -18245.8, -16789)
    ,@("LTW", 132, 6764.1274, 6130.048, 8812.691999999999, 18390.144, 26438.076, -15860.144, -31498.076)
    ,@("WVR", 132, 1331.0358, 1279.6154, 1999.5, 3838.8462, 5998.5, -1308.8462, -11058.5)
    ,@("WVR", 136, 3655.2896, 2677.6897, 6805.3335, 8033.0691, 20416.0005, -5483.0691, -25516.0005)
)

$currentSheetName = $null
$ws = $null

foreach ($edit in $edits) {
    $sheetName = $edit[0]
    $row = $edit[1]

    if ($sheetName -ne $currentSheetName) {
        $ws = $wb.Worksheets.Item($sheetName)
        $currentSheetName = $sheetName
    }

    for ($i = 0; $i -lt 7; $i++) {
        $value = $edit[2 + $i]
        if ($null -eq $value) {
            continue
        }
        $col = 8 + $i  # H=8 .. N=14
        $cell = $ws.Cells.Item($row, $col)
        if ("$value" -eq $CLEAR) {
            $cell.Value = ""
        } else {
            $cell.Value = $value
        }
    }
}
